$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-detected as a number by Excel;
# force them to stay text (matches the source workbook, where every Price/Volume
# cell is stored as a text string) by pre-setting the number format to "@".
$textCells = @("D4","D5","D6","D10","D12","D13","D19","D20","D22","D25","D34","D35","D36","D39","D40","D43","D44","D46","D48","D49","D50")
foreach ($addr in $textCells) { $ws.Range($addr).NumberFormat = "@" }

# Apply the updated coin price / volume figures (and the swapped Filecoin /
# dogwifhat rows 39-40) from the latest GitHub Actions refresh.
$ws.Range("D2").Value = '67.756.41'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '3.794.93'
$ws.Range("E3").Value = '  +0.39%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.18%  '
$ws.Range("D5").Value = '600.09'
$ws.Range("E5").Value = '  +0.78%  '
$ws.Range("D6").Value = '165.13'
$ws.Range("E6").Value = '  -1.07%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("E8").Value = '  -0.48%  '
$ws.Range("E9").Value = '  -0.19%  '
$ws.Range("D10").Value = '0.452'
$ws.Range("E10").Value = '  +0.58%  '
$ws.Range("E11").Value = '  +1.44%  '
$ws.Range("D12").Value = '0.0000249'
$ws.Range("E12").Value = '  -1.33%  '
$ws.Range("D13").Value = '35.77'
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("D14").Value = '4.435.88'
$ws.Range("E14").Value = '  +0.48%  '
$ws.Range("D15").Value = '3.801.67'
$ws.Range("E15").Value = '  +1.19%  '
$ws.Range("D16").Value = '67.755.00'
$ws.Range("E16").Value = '  +0.22%  '
$ws.Range("E17").Value = '  -0.71%  '
$ws.Range("E18").Value = '  +1.27%  '
$ws.Range("D19").Value = '7.05'
$ws.Range("E19").Value = '  +0.48%  '
$ws.Range("D20").Value = '463.45'
$ws.Range("E20").Value = '  +0.82%  '
$ws.Range("E21").Value = '  -2.55%  '
$ws.Range("D22").Value = '0.700'
$ws.Range("E22").Value = '  +0.51%  '
$ws.Range("E23").Value = '  -5.26%  '
$ws.Range("E24").Value = '  -0.49%  '
$ws.Range("D25").Value = '12.04'
$ws.Range("E25").Value = '  +0.50%  '
$ws.Range("E26").Value = '  -0.64%  '
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("E28").Value = '  -0.06%  '
$ws.Range("D29").Value = '3.945.30'
$ws.Range("E29").Value = '  +0.45%  '
$ws.Range("E30").Value = '  -0.11%  '
$ws.Range("E31").Value = '  +2.45%  '
$ws.Range("E33").Value = '  -1.20%  '
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("D35").Value = '9.02'
$ws.Range("E35").Value = '  -0.62%  '
$ws.Range("D36").Value = '0.0995'
$ws.Range("E36").Value = '  -0.58%  '
$ws.Range("E37").Value = '  +0.47%  '
$ws.Range("E38").Value = '  +0.07%  '
$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D39").Value = '5.76'
$ws.Range("E39").Value = '  -0.15%  '
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").Value = '3.22'
$ws.Range("E40").Value = '  -4.38%  '
$ws.Range("E41").Value = '  +0.10%  '
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("D43").Value = '45.28'
$ws.Range("E43").Value = '  -0.38%  '
$ws.Range("D44").Value = '47.63'
$ws.Range("E44").Value = '  -0.89%  '
$ws.Range("E45").Value = '  -0.43%  '
$ws.Range("D46").Value = '151.08'
$ws.Range("E46").Value = '  +0.89%  '
$ws.Range("E47").Value = '  +10.61%  '
$ws.Range("D48").Value = '27.62'
$ws.Range("E48").Value = '  +3.40%  '
$ws.Range("D49").Value = '8.33'
$ws.Range("E49").Value = '  +0.29%  '
$ws.Range("D50").Value = '393.88'
$ws.Range("E50").Value = '  -0.02%  '
$ws.Range("E51").Value = '  +1.63%  '
